# WDZO-1321: changed templates for uploading/deleting users
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update the "Точка продажи" sample code ---
$ws.Range("E2").Value = "038_9038_0393"

# --- Row 3: second sample employee -> д'Артаньян / Артем / Петрович ---
$ws.Range("B3").Value = "д’Артаньян"
$ws.Range("C3").Value = "Артем"
$ws.Range("D3").Value = "Петрович"
$ws.Range("E3").Value = "038_9038_0393"

# --- Row 4: third sample employee -> Плотникова-Работникова / Екатерина / Федоровна ---
$ws.Range("B4").Value = "Плотникова-Работникова"
$ws.Range("C4").Value = "Екатерина"
$ws.Range("D4").Value = "Федоровна"
$ws.Range("E4").Value = "038_9038_0393"

# Give the re-typed example name cells (B3:D4) their own (new) plain font entry,
# matching the font fork Excel produced when the sample rows were edited.
$ws.Range("B3:D4").Font.Size = 12.5

# --- Rows 5-7: blank template rows (clear sample data, keep formatting) ---
$ws.Range("A5:E7").ClearContents()

# --- Row 8: instructional note in column A, rest cleared ---
$ws.Range("A8").Value = "Удали содержимое и используй как шаблон ツ"
$ws.Range("B8:E8").ClearContents()

# --- Row 9: fully cleared; A9/B9 highlighted in red ---
$ws.Range("A9:E9").ClearContents()
$ws.Range("A9:B9").Font.Color = 192

# --- Rows 10-11: blank template rows ---
$ws.Range("A10:E11").ClearContents()

# --- Selection left where the last edit happened ---
$ws.Range("B16").Select()
